$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (current rows 459-460),
# pushing all existing data rows down by two.
$ws.Rows("459:460").Insert()

# New row 459 data
$ws.Cells.Item(459, 1).Value = 9
$ws.Cells.Item(459, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(459, 3).Value = "Metropolitana"
$ws.Cells.Item(459, 4).Value = 44946
$ws.Cells.Item(459, 5).Value = 13
$ws.Cells.Item(459, 6).Value = 100112032
$ws.Cells.Item(459, 7).Value = "Zapallo italiano"
$ws.Cells.Item(459, 8).Value = "Sin especificar"
$ws.Cells.Item(459, 9).Value = "Primera"
$ws.Cells.Item(459, 10).Value = 430
$ws.Cells.Item(459, 11).Value = 6000
$ws.Cells.Item(459, 12).Value = 7000
$ws.Cells.Item(459, 13).Value = 6500
$ws.Cells.Item(459, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(459, 15).Value = "Región Metropolitana"
$ws.Cells.Item(459, 16).Value = 130
$ws.Cells.Item(459, 17).Value = 50
$ws.Cells.Item(459, 18).Value = "Hortaliza"

# New row 460 data
$ws.Cells.Item(460, 1).Value = 9
$ws.Cells.Item(460, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(460, 3).Value = "Metropolitana"
$ws.Cells.Item(460, 4).Value = 44946
$ws.Cells.Item(460, 5).Value = 13
$ws.Cells.Item(460, 6).Value = 100112032
$ws.Cells.Item(460, 7).Value = "Zapallo italiano"
$ws.Cells.Item(460, 8).Value = "Sin especificar"
$ws.Cells.Item(460, 9).Value = "Primera"
$ws.Cells.Item(460, 10).Value = 340
$ws.Cells.Item(460, 11).Value = 6000
$ws.Cells.Item(460, 12).Value = 7000
$ws.Cells.Item(460, 13).Value = 6500
$ws.Cells.Item(460, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(460, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(460, 16).Value = 130
$ws.Cells.Item(460, 17).Value = 50
$ws.Cells.Item(460, 18).Value = "Hortaliza"

Write-Output "Done"
